$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.851.40'
$ws.Range("E2").Value = '  -1.19%  '

$ws.Range("D3").Value = '1.890.65'
$ws.Range("E3").Value = '  -1.38%  '

$ws.Range("E4").Value = '  -0.14%  '

$c = $ws.Range("D5")
$c.Value = "'0.7754"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -3.26%  '

$c = $ws.Range("D6")
$c.Value = "'244.69"
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.15%  '

$ws.Range("E7").Value = '  -0.12%  '

$c = $ws.Range("D8")
$c.Value = "'0.3144"
$c.Style = "Normal"
$ws.Range("E8").Value = '  -3.08%  '

$c = $ws.Range("D9")
$c.Value = "'0.07407"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +2.64%  '

$c = $ws.Range("D10")
$c.Value = "'25.34"
$c.Style = "Normal"
$ws.Range("E10").Value = '  -5.49%  '

$c = $ws.Range("D11")
$c.Value = "'0.08129"
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.65%  '

$c = $ws.Range("D12")
$c.Value = "'0.7675"
$c.Style = "Normal"
$ws.Range("E12").Value = '  -2.29%  '

$c = $ws.Range("D13")
$c.Value = "'5.475"
$c.Style = "Normal"
$ws.Range("E13").Value = '  +1.56%  '

$ws.Range("D14").Value = '1.885.42'
$ws.Range("E14").Value = '  -1.47%  '

$c = $ws.Range("D15")
$c.Value = "'92.30"
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.45%  '

$c = $ws.Range("D16")
$c.Value = "'6.174"
$c.Style = "Normal"
$ws.Range("E16").Value = '  +2.18%  '

$ws.Range("D17").Value = '29.887.58'
$ws.Range("E17").Value = '  -1.16%  '

$c = $ws.Range("D18")
$c.Value = "'13.97"
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.59%  '

$c = $ws.Range("D19")
$c.Value = "'244.67"
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.85%  '

$c = $ws.Range("D20")
$c.Value = "'0.000007851"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.27%  '

$ws.Range("E21").Value = '  -0.11%  '

$c = $ws.Range("D22")
$c.Value = "'8.100"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.91%  '

$ws.Range("D23").Value = '2.129.13'
$ws.Range("E23").Value = '  -2.03%  '

$c = $ws.Range("D24")
$c.Value = "'1.0000"
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.34%  '

$c = $ws.Range("D25")
$c.Value = "'0.1582"
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.87%  '

$c = $ws.Range("D26")
$c.Value = "'9.420"
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.32%  '

$c = $ws.Range("D27")
$c.Value = "'162.53"
$c.Style = "Normal"
$ws.Range("E27").Value = '  -2.85%  '

$c = $ws.Range("D28")
$c.Value = "'18.81"
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.86%  '

$c = $ws.Range("D29")
$c.Value = "'2.039"
$c.Style = "Normal"
$ws.Range("E29").Value = '  -4.85%  '

$c = $ws.Range("D30")
$c.Value = "'1.444"
$c.Style = "Normal"
$ws.Range("E30").Value = '  +3.96%  '

$ws.Range("E31").Value = '  -0.06%  '

$c = $ws.Range("D32")
$c.Value = "'4.501"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.26%  '

$ws.Range("E33").Value = '  -1.70%  '

$c = $ws.Range("D34")
$c.Value = "'0.05556"
$c.Style = "Normal"
$ws.Range("E34").Value = '  -2.26%  '

$c = $ws.Range("D35")
$c.Value = "'1.246"
$c.Style = "Normal"
$ws.Range("E35").Value = '  -3.48%  '

$c = $ws.Range("D36")
$c.Value = "'0.7604"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.55%  '

$c = $ws.Range("D37")
$c.Value = "'1.002"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.09%  '

$c = $ws.Range("D38")
$c.Value = "'2.644"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -3.25%  '

$c = $ws.Range("D39")
$c.Value = "'0.01932"
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.42%  '

$c = $ws.Range("D40")
$c.Value = "'2.787"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.30%  '

$ws.Range("D41").Value = '1.161.98'
$ws.Range("E41").Value = '  +11.95%  '

$c = $ws.Range("D42")
$c.Value = "'0.4464"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.08%  '

$c = $ws.Range("D43")
$c.Value = "'73.98"
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.36%  '

$c = $ws.Range("D44")
$c.Value = "'5.970"
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.74%  '

$c = $ws.Range("D45")
$c.Value = "'0.8490"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.79%  '

$c = $ws.Range("D46")
$c.Value = "'1.000"
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.14%  '

$c = $ws.Range("D47")
$c.Value = "'1.899"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.23%  '

$c = $ws.Range("D48")
$c.Value = "'102.24"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.79%  '

$c = $ws.Range("D49")
$c.Value = "'9.932"
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.70%  '

$c = $ws.Range("D50")
$c.Value = "'3.080"
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.60%  '

$c = $ws.Range("D51")
$c.Value = "'7.524"
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.16%  '
